# Horarios actualizados Línea 141 - 1283
# Applies the scrape refresh (new "Última actualización" / "Total filas" +
# updated / inserted schedule rows) to all three worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 11:01:19"
$ws1.Range("A3").Value = "Total filas: 121"

function Set-Row($ws, $r, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

# Rows whose scrape-time / arrival swapped with a neighbour
Set-Row $ws1 55  "08:14:55" "08:53" "215B_EL PATO"       39  "LP1912"
Set-Row $ws1 56  "08:49:06" "08:53" "23_HERNANDEZ"        4  "LP1912"
Set-Row $ws1 71  "08:49:06" "09:31" "23_HERNANDEZ"       42  "LP1912"
Set-Row $ws1 72  "08:14:55" "09:31" "16_SANTA ANA"       77  "LP1912"

# "14_ABASTO" / "15_ABASTO" line-name swap
$ws1.Cells.Item(90, 3).Value = "15_ABASTO"
$ws1.Cells.Item(91, 3).Value = "14_ABASTO"

# Rows 99-115: refreshed scrape time / updated minutes, row 105/111/114 keep
# their original "Hora_Scrap" (still from the prior scrape cycle)
Set-Row $ws1 99  "11:01:19" "11:02" "81_EL PELIGRO"        1 "LP1912"
Set-Row $ws1 100 "11:01:19" "11:06" "23_HERNANDEZ"         5 "LP1912"
Set-Row $ws1 101 "11:01:19" "11:10" "16_P MOR-SANTA ANA"   9 "LP1912"
Set-Row $ws1 102 "11:01:19" "11:14" "14_ABASTO"           13 "LP1912"
Set-Row $ws1 104 "11:01:19" "11:15" "15X38_ABASTO"        14 "LP1912"
Set-Row $ws1 105 "10:32:07" "11:24" "16_SANTA ANA"        52 "LP1912"
Set-Row $ws1 106 "11:01:19" "11:25" "16_SANTA ANA"        24 "LP1912"
Set-Row $ws1 107 "11:01:19" "11:29" "10_OLMOS"             28 "LP1912"
Set-Row $ws1 108 "09:42:42" "11:30" "215C_EL PATO"        108 "LP1912"
Set-Row $ws1 109 "11:01:19" "11:31" "215C_EL PATO"         30 "LP1912"
Set-Row $ws1 110 "11:01:19" "11:41" "215B_EL PATO"         40 "LP1912"
Set-Row $ws1 111 "10:32:07" "11:42" "215B_EL PATO"         70 "LP1912"
Set-Row $ws1 112 "11:01:19" "11:45" "15X38_ABASTO"         44 "LP1912"
Set-Row $ws1 113 "11:01:19" "11:47" "23_HERNANDEZ"         46 "LP1912"
Set-Row $ws1 114 "10:32:07" "11:51" "23_HERNANDEZ"         79 "LP1912"
Set-Row $ws1 115 "11:01:19" "11:53" "225_GOMEZ"            52 "LP1912"

# New rows appended by the refreshed scrape (116-126)
Set-Row $ws1 116 "11:01:19" "11:58" "17_ROMERO"            57 "LP1912"
Set-Row $ws1 117 "11:01:19" "12:05" "11_ETCHEVERRY"        64 "LP1912"
Set-Row $ws1 118 "10:32:07" "12:06" "11_ETCHEVERRY"        94 "LP1912"
Set-Row $ws1 119 "11:01:19" "12:10" "15_ABASTO"            69 "LP1912"
Set-Row $ws1 120 "11:01:19" "12:10" "16_P MOR-SANTA ANA"   69 "LP1912"
Set-Row $ws1 121 "11:01:19" "12:17" "10_OLMOS"             76 "LP1912"
Set-Row $ws1 122 "11:01:19" "12:22" "215C_EL PATO"         81 "LP1912"
Set-Row $ws1 123 "11:01:19" "12:32" "14_ABASTO"            91 "LP1912"
Set-Row $ws1 124 "11:01:19" "12:34" "15_ABASTO"            93 "LP1912"
Set-Row $ws1 125 "11:01:19" "12:37" "27_EL RETIRO"         96 "LP1912"
Set-Row $ws1 126 "11:01:19" "12:48" "16_SANTA ANA"        107 "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 11:01:19"
$ws2.Range("A3").Value = "Total filas: 21"

Set-Row $ws2 23 "11:01:19" "11:31" "215C_EL PATO" 30 "LP1912"
Set-Row $ws2 24 "11:01:19" "11:41" "215B_EL PATO" 40 "LP1912"
Set-Row $ws2 25 "10:32:07" "11:42" "215B_EL PATO" 70 "LP1912"
Set-Row $ws2 26 "11:01:19" "12:22" "215C_EL PATO" 81 "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 11:01:19"

Set-Row $ws3 25 "11:01:19" "11:26" "215C_LA PLATA" 25 "L6203"
